$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 837.71875
$ws.Range("I15").Value = 837.71875
$ws.Range("K15").Value = 2513.15625
$ws.Range("M15").Value = -2344.15625

$ws.Range("H28").Value = 2481.2222
$ws.Range("I28").Value = 314
$ws.Range("K28").Value = 314
$ws.Range("M28").Value = 171

$ws.Range("H101").Value = 703.125
$ws.Range("I101").Value = 236.8
$ws.Range("J101").Value = 1480.3334
$ws.Range("K101").Value = 710.4000000000001
$ws.Range("L101").Value = 4441.0002
$ws.Range("M101").Value = 911.5999999999999
$ws.Range("N101").Value = -7685.0002

$ws.Range("H113").Value = 4681.3335
$ws.Range("I113").Value = 4225.8823
$ws.Range("J113").Value = 5787.4287
$ws.Range("K113").Value = 4225.8823
$ws.Range("L113").Value = 5787.4287
$ws.Range("M113").Value = -971.8823000000002
$ws.Range("N113").Value = -12295.4287

$ws.Range("H116").Value = 10749.833
$ws.Range("I116").Value = 12124.75
$ws.Range("J116").Value = 8000
$ws.Range("K116").Value = 12124.75
$ws.Range("L116").Value = 8000
$ws.Range("M116").Value = -8682.75
$ws.Range("N116").Value = -14884

$ws.Range("H132").Value = 28574120
$ws.Range("I132").Value = 32260938
$ws.Range("J132").Value = 1292.5
$ws.Range("K132").Value = 96782814
$ws.Range("L132").Value = 3877.5
$ws.Range("M132").Value = -96780284
$ws.Range("N132").Value = -8937.5

$ws.Range("H137").Value = 6872.9614
$ws.Range("I137").Value = 2261.9333
$ws.Range("J137").Value = 13160.728
$ws.Range("K137").Value = 6785.7999
$ws.Range("L137").Value = 39482.18399999999
$ws.Range("M137").Value = -4235.7999
$ws.Range("N137").Value = -44582.18399999999

$ws.Range("H140").Value = 139999.77
$ws.Range("J140").Value = 139999.77
$ws.Range("L140").Value = 139999.77
$ws.Range("N140").Value = -150359.77

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2506
$ws.Range("I2").Value = 1945.9524
$ws.Range("J2").Value = 4466.1665
$ws.Range("K2").Value = 1945.9524
$ws.Range("L2").Value = 4466.1665
$ws.Range("M2").Value = -1832.9524
$ws.Range("N2").Value = -4692.1665

$ws.Range("H45").Value = 4229.8
$ws.Range("I45").Value = 4599.6665
$ws.Range("K45").Value = 4599.6665
$ws.Range("M45").Value = -4222.6665

$ws.Range("H63").Value = 9389.444
$ws.Range("I63").Value = 8168.3335
$ws.Range("K63").Value = 8168.3335
$ws.Range("M63").Value = -7482.3335

$ws.Range("H66").Value = 9389.444
$ws.Range("I66").Value = 8168.3335
$ws.Range("K66").Value = 40841.6675
$ws.Range("M66").Value = -37409.6675

$ws.Range("H97").Value = 2211.65
$ws.Range("I97").Value = 1639.75
$ws.Range("K97").Value = 1639.75
$ws.Range("M97").Value = -1143.75

$ws.Range("H110").Value = 3658.5217
$ws.Range("I110").Value = 4942
$ws.Range("K110").Value = 4942
$ws.Range("M110").Value = -2897

$ws.Range("H116").Value = 2506
$ws.Range("I116").Value = 1945.9524
$ws.Range("J116").Value = 4466.1665
$ws.Range("K116").Value = 1945.9524
$ws.Range("L116").Value = 4466.1665
$ws.Range("M116").Value = 348.0476000000001
$ws.Range("N116").Value = -9054.166499999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2506
$ws.Range("I3").Value = 1945.9524
$ws.Range("J3").Value = 4466.1665
$ws.Range("K3").Value = 1945.9524
$ws.Range("L3").Value = 4466.1665
$ws.Range("M3").Value = -1831.9524
$ws.Range("N3").Value = -4694.1665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 2108589.5
$ws.Range("I4").Value = 4004720
$ws.Range("J4").Value = 1777.7778
$ws.Range("K4").Value = 4004720
$ws.Range("L4").Value = 1777.7778
$ws.Range("M4").Value = -4004608
$ws.Range("N4").Value = -2001.7778

$ws.Range("H31").Value = 43108010
$ws.Range("J31").Value = 156259500
$ws.Range("L31").Value = 156259500
$ws.Range("N31").Value = -156260090

$ws.Range("H34").Value = 43108010
$ws.Range("J34").Value = 156259500
$ws.Range("L34").Value = 156259500
$ws.Range("N34").Value = -156259904

$ws.Range("H62").Value = 3178.0557
$ws.Range("I62").Value = 3206.0588
$ws.Range("K62").Value = 3206.0588
$ws.Range("M62").Value = -2582.0588

$ws.Range("H65").Value = 3178.0557
$ws.Range("I65").Value = 3206.0588
$ws.Range("K65").Value = 16030.294
$ws.Range("M65").Value = -12910.294

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 467.9
$ws.Range("I92").Value = 322.5
$ws.Range("J92").Value = 564.8333
$ws.Range("K92").Value = 967.5
$ws.Range("L92").Value = 1694.4999
$ws.Range("M92").Value = 280.5
$ws.Range("N92").Value = -4190.4999

$ws.Range("H113").Value = 1921.3334
$ws.Range("I113").Value = 1632.125
$ws.Range("K113").Value = 4896.375
$ws.Range("M113").Value = -2726.375

$ws.Range("H122").Value = 1847.1428
$ws.Range("I122").Value = 825
$ws.Range("J122").Value = 2256
$ws.Range("K122").Value = 7425
$ws.Range("L122").Value = 20304
$ws.Range("M122").Value = -4975
$ws.Range("N122").Value = -25204

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1724353.1
$ws.Range("I2").Value = 3125093.5
$ws.Range("J2").Value = 365.07693
$ws.Range("K2").Value = 3125093.5
$ws.Range("L2").Value = 365.07693
$ws.Range("M2").Value = -3124980.5
$ws.Range("N2").Value = -591.0769299999999

$ws.Range("H18").Value = 29999.334
$ws.Range("I18").Value = 29999.334
$ws.Range("K18").Value = 29999.334
$ws.Range("M18").Value = -29706.334

$ws.Range("H43").Value = 3432.4285
$ws.Range("I43").Value = 3432.4285
$ws.Range("K43").Value = 3432.4285
$ws.Range("M43").Value = -3281.4285

$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("L57").ClearContents()
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = 0

$ws.Range("H80").Value = 6649
$ws.Range("I80").Value = 7983.1665
$ws.Range("K80").Value = 7983.1665
$ws.Range("M80").Value = -6985.1665

$ws.Range("H83").Value = 6649
$ws.Range("I83").Value = 7983.1665
$ws.Range("K83").Value = 39915.8325
$ws.Range("M83").Value = -34923.8325

$ws.Range("H97").Value = 696.04346
$ws.Range("I97").Value = 646.3333
$ws.Range("J97").Value = 789.25
$ws.Range("K97").Value = 646.3333
$ws.Range("L97").Value = 789.25
$ws.Range("M97").Value = -150.3333
$ws.Range("N97").Value = -1781.25

$ws.Range("H113").Value = 3355.4473
$ws.Range("I113").Value = 2851.1
$ws.Range("K113").Value = 2851.1
$ws.Range("M113").Value = -681.0999999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1992.5405
$ws.Range("J46").Value = 4797.5
$ws.Range("L46").Value = 4797.5
$ws.Range("N46").Value = -5173.5

$ws.Range("H55").Value = 660.2632
$ws.Range("I55").Value = 534
$ws.Range("K55").Value = 534
$ws.Range("M55").Value = -361

$ws.Range("H68").Value = 2684.4546
$ws.Range("I68").Value = 1984.25
$ws.Range("K68").Value = 1984.25
$ws.Range("M68").Value = -1235.25

$ws.Range("H71").Value = 2684.4546
$ws.Range("I71").Value = 1984.25
$ws.Range("K71").Value = 9921.25
$ws.Range("M71").Value = -6177.25

$ws.Range("H100").Value = 5024.7144
$ws.Range("I100").Value = 3957.889
$ws.Range("J100").Value = 6945
$ws.Range("K100").Value = 3957.889
$ws.Range("L100").Value = 6945
$ws.Range("M100").Value = -3416.889
$ws.Range("N100").Value = -8027

$ws.Range("H122").Value = 4171242.5
$ws.Range("I122").Value = 4141.8096
$ws.Range("K122").Value = 12425.4288
$ws.Range("M122").Value = -9975.428799999998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5551.4375
$ws.Range("I81").Value = 4649.8
$ws.Range("J81").Value = 5961.273
$ws.Range("K81").Value = 9299.6
$ws.Range("L81").Value = 11922.546
$ws.Range("M81").Value = -8238.6
$ws.Range("N81").Value = -14044.546

$ws.Range("H84").Value = 5551.4375
$ws.Range("I84").Value = 4649.8
$ws.Range("J84").Value = 5961.273
$ws.Range("K84").Value = 46498
$ws.Range("L84").Value = 59612.73
$ws.Range("M84").Value = -41194
$ws.Range("N84").Value = -70220.73000000001
